# Updates the cryptos price table (Price / Volume(1h) columns, and two
# coin-row swaps) to match the refreshed data feed from the commit:
# "Updated cryptos list on Mon Jan  1 23:46:14 UTC 2024 with GitHub Actions"
#
# Every written value is prefixed with a leading apostrophe (quote-prefix)
# so Excel stores it as literal text -- matching the workbook's existing
# inline-string cells -- instead of auto-coercing numeric-looking strings
# (e.g. "7.60", "2.30") into numbers and silently dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'44.084.77"
$ws.Cells.Item(2, 5).Value = "'  +4.37%  "

$ws.Cells.Item(3, 4).Value = "'2.348.96"
$ws.Cells.Item(3, 5).Value = "'  +3.04%  "

$ws.Cells.Item(4, 5).Value = "'  +0.36%  "

$ws.Cells.Item(5, 4).Value = "'315.44"
$ws.Cells.Item(5, 5).Value = "'  +0.91%  "

$ws.Cells.Item(6, 4).Value = "'109.75"
$ws.Cells.Item(6, 5).Value = "'  +7.83%  "

$ws.Cells.Item(7, 5).Value = "'  +3.57%  "

$ws.Cells.Item(8, 5).Value = "'  +0.28%  "

$ws.Cells.Item(9, 5).Value = "'  +5.13%  "

$ws.Cells.Item(10, 4).Value = "'42.06"
$ws.Cells.Item(10, 5).Value = "'  +8.61%  "

$ws.Cells.Item(11, 4).Value = "'0.0923"
$ws.Cells.Item(11, 5).Value = "'  +2.91%  "

$ws.Cells.Item(12, 5).Value = "'  +5.54%  "

$ws.Cells.Item(13, 5).Value = "'  +4.80%  "

$ws.Cells.Item(14, 5).Value = "'  +0.32%  "

$ws.Cells.Item(15, 4).Value = "'15.58"
$ws.Cells.Item(15, 5).Value = "'  +4.24%  "

$ws.Cells.Item(16, 4).Value = "'2.703.17"
$ws.Cells.Item(16, 5).Value = "'  +2.96%  "

$ws.Cells.Item(17, 4).Value = "'2.350.27"
$ws.Cells.Item(17, 5).Value = "'  +3.20%  "

$ws.Cells.Item(18, 4).Value = "'44.285.02"
$ws.Cells.Item(18, 5).Value = "'  +4.94%  "

$ws.Cells.Item(19, 4).Value = "'7.60"
$ws.Cells.Item(19, 5).Value = "'  +5.23%  "

$ws.Cells.Item(20, 5).Value = "'  +2.95%  "

$ws.Cells.Item(21, 4).Value = "'13.04"
$ws.Cells.Item(21, 5).Value = "'  -2.97%  "

$ws.Cells.Item(22, 4).Value = "'74.71"
$ws.Cells.Item(22, 5).Value = "'  +2.63%  "

$ws.Cells.Item(23, 5).Value = "'  +0.45%  "

$ws.Cells.Item(24, 4).Value = "'269.48"
$ws.Cells.Item(24, 5).Value = "'  +2.30%  "

$ws.Cells.Item(25, 4).Value = "'2.30"
$ws.Cells.Item(25, 5).Value = "'  +6.70%  "

$ws.Cells.Item(26, 5).Value = "'  -0.19%  "

$ws.Cells.Item(27, 4).Value = "'7.68"
$ws.Cells.Item(27, 5).Value = "'  +11.45%  "

$ws.Cells.Item(28, 4).Value = "'11.23"
$ws.Cells.Item(28, 5).Value = "'  +5.80%  "

$ws.Cells.Item(29, 5).Value = "'  +1.90%  "

$ws.Cells.Item(30, 4).Value = "'39.71"
$ws.Cells.Item(30, 5).Value = "'  +10.74%  "

$ws.Cells.Item(31, 4).Value = "'22.74"
$ws.Cells.Item(31, 5).Value = "'  +2.73%  "

$ws.Cells.Item(32, 2).Value = "'Hedera"
$ws.Cells.Item(32, 3).Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).Value = "'0.0926"
$ws.Cells.Item(32, 5).Value = "'  +7.96%  "

$ws.Cells.Item(33, 2).Value = "'Monero"
$ws.Cells.Item(33, 3).Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(33, 4).Value = "'169.23"
$ws.Cells.Item(33, 5).Value = "'  +3.20%  "

$ws.Cells.Item(34, 4).Value = "'2.88"
$ws.Cells.Item(34, 5).Value = "'  +10.53%  "

$ws.Cells.Item(35, 5).Value = "'  +2.01%  "

$ws.Cells.Item(36, 5).Value = "'  +4.67%  "

$ws.Cells.Item(37, 4).Value = "'4.74"
$ws.Cells.Item(37, 5).Value = "'  +6.23%  "

$ws.Cells.Item(38, 4).Value = "'0.0367"
$ws.Cells.Item(38, 5).Value = "'  +6.19%  "

$ws.Cells.Item(39, 5).Value = "'  +10.50%  "

$ws.Cells.Item(40, 4).Value = "'3.82"
$ws.Cells.Item(40, 5).Value = "'  +3.98%  "

$ws.Cells.Item(41, 4).Value = "'1.73"
$ws.Cells.Item(41, 5).Value = "'  +10.52%  "

$ws.Cells.Item(42, 4).Value = "'105.86"
$ws.Cells.Item(42, 5).Value = "'  +8.27%  "

$ws.Cells.Item(43, 4).Value = "'13.94"
$ws.Cells.Item(43, 5).Value = "'  +17.52%  "

$ws.Cells.Item(44, 4).Value = "'0.242"
$ws.Cells.Item(44, 5).Value = "'  +6.97%  "

$ws.Cells.Item(45, 4).Value = "'72.02"
$ws.Cells.Item(45, 5).Value = "'  +4.69%  "

$ws.Cells.Item(46, 5).Value = "'  +0.54%  "

$ws.Cells.Item(47, 4).Value = "'115.98"
$ws.Cells.Item(47, 5).Value = "'  +5.70%  "

$ws.Cells.Item(48, 4).Value = "'78.14"
$ws.Cells.Item(48, 5).Value = "'  -0.84%  "

$ws.Cells.Item(49, 4).Value = "'1.661.74"
$ws.Cells.Item(49, 5).Value = "'  -2.20%  "

$ws.Cells.Item(50, 2).Value = "'FraxShare"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(50, 4).Value = "'9.03"
$ws.Cells.Item(50, 5).Value = "'  +4.72%  "

$ws.Cells.Item(51, 2).Value = "'TheGraph"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(51, 4).Value = "'0.218"
$ws.Cells.Item(51, 5).Value = "'  +17.69%  "

